# The edit reshuffles the 20 data rows (rows 2-21) of the single sheet:
# each destination row ends up holding the exact data that used to live in
# a (different) source row. Row 1 (headers) is untouched. Columns span
# A..R (18 columns).
#
# Mapping: destination row -> source row (both are pre-edit row numbers)
$rowMap = @{
    2  = 3
    3  = 4
    4  = 6
    5  = 18
    6  = 19
    7  = 20
    8  = 10
    9  = 11
    10 = 12
    11 = 7
    12 = 8
    13 = 15
    14 = 13
    15 = 21
    16 = 14
    17 = 9
    18 = 5
    19 = 16
    20 = 17
    21 = 2
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1
$lastCol = 18   # columns A (1) .. R (18)

# 1) Snapshot every source row's values (row number -> array of column values)
#    before any writes happen, since several rows are both sources and
#    destinations for other rows.
$snapshot = @{}
for ($r = 2; $r -le 21; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each destination row from the snapshot of its mapped source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c - 1]
    }
}
